$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K", formerly Strike#) values per recalculated s_vals
$kValues = @{
    2 = 4
    3 = 0
    4 = 0
    5 = 2
    6 = 2
    7 = 2
    8 = 2
    9 = 1
    10 = 3
    11 = 2
    12 = 2
    13 = 1
    14 = 2
    15 = 1
    16 = 1
    17 = 0
    18 = 1
    19 = 1
    20 = 2
    21 = 1
    22 = 1
    23 = 4
    24 = 1
    25 = 0
    26 = 0
    27 = 1
    28 = 1
    29 = 0
    30 = 1
    31 = 2
    32 = 4
    35 = 1
    36 = 1
    37 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

